$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2311.8333
$ws.Range("I19").Value = 2242.75
$ws.Range("K19").Value = 2242.75
$ws.Range("M19").Value = -2067.75
$ws.Range("H28").Value = 2824.25
$ws.Range("I28").Value = 432.33334
$ws.Range("K28").Value = 432.33334
$ws.Range("M28").Value = 52.66665999999998
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H107").Value = 1996
$ws.Range("I107").Value = 648
$ws.Range("J107").Value = 3044.4443
$ws.Range("K107").Value = 648
$ws.Range("L107").Value = 3044.4443
$ws.Range("M107").Value = 1272
$ws.Range("N107").Value = -6884.4443
$ws.Range("H116").Value = 5233.4
$ws.Range("I116").Value = 4822.3335
$ws.Range("K116").Value = 4822.3335
$ws.Range("M116").Value = -1380.3335
$ws.Range("H137").Value = 2399.6667
$ws.Range("I137").Value = 1999
$ws.Range("K137").Value = 5997
$ws.Range("M137").Value = -3447
$ws.Range("H138").Value = 2548.7058
$ws.Range("I138").Value = 761.1429000000001
$ws.Range("J138").Value = 3800
$ws.Range("K138").Value = 2283.4287
$ws.Range("L138").Value = 11400
$ws.Range("M138").Value = 2856.5713
$ws.Range("N138").Value = -21680

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2174.4
$ws.Range("I2").Value = 1278.2307
$ws.Range("K2").Value = 1278.2307
$ws.Range("M2").Value = -1165.2307
$ws.Range("H45").Value = 1849.5454
$ws.Range("J45").Value = 2499
$ws.Range("L45").Value = 2499
$ws.Range("N45").Value = -3253
$ws.Range("H61").Value = 3611.625
$ws.Range("I61").Value = 3474
$ws.Range("K61").Value = 3474
$ws.Range("M61").Value = -3262
$ws.Range("H74").Value = 3075
$ws.Range("I74").Value = 2365
$ws.Range("K74").Value = 2365
$ws.Range("M74").Value = -1491
$ws.Range("H77").Value = 3075
$ws.Range("I77").Value = 2365
$ws.Range("K77").Value = 11825
$ws.Range("M77").Value = -7457
$ws.Range("H97").Value = 871
$ws.Range("I97").Value = 896.5454999999999
$ws.Range("K97").Value = 896.5454999999999
$ws.Range("M97").Value = -400.5454999999999
$ws.Range("H116").Value = 2174.4
$ws.Range("I116").Value = 1278.2307
$ws.Range("K116").Value = 1278.2307
$ws.Range("M116").Value = 1015.7693
$ws.Range("H132").Value = 1082.3334
$ws.Range("I132").Value = 1082.3334
$ws.Range("K132").Value = 3247.0002
$ws.Range("M132").Value = -717.0001999999999
$ws.Range("H136").Value = 3611.625
$ws.Range("I136").Value = 3474
$ws.Range("K136").Value = 10422
$ws.Range("M136").Value = -7872

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2174.4
$ws.Range("I3").Value = 1278.2307
$ws.Range("K3").Value = 1278.2307
$ws.Range("M3").Value = -1164.2307
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 2792.4285
$ws.Range("I94").Value = 3099.4546
$ws.Range("J94").Value = 1666.6666
$ws.Range("K94").Value = 3099.4546
$ws.Range("L94").Value = 1666.6666
$ws.Range("M94").Value = -2648.4546
$ws.Range("N94").Value = -2568.6666
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3018.75
$ws.Range("I58").Value = 2664.2856
$ws.Range("K58").Value = 2664.2856
$ws.Range("M58").Value = -2461.2856
$ws.Range("H132").Value = 735.6667
$ws.Range("I132").Value = 735.6667
$ws.Range("K132").Value = 2207.0001
$ws.Range("M132").Value = 322.9998999999998
$ws.Range("H134").Value = 2093.4167
$ws.Range("I134").Value = 1966.4348
$ws.Range("K134").Value = 5899.3044
$ws.Range("M134").Value = -3364.3044
$ws.Range("H136").Value = 3018.75
$ws.Range("I136").Value = 2664.2856
$ws.Range("K136").Value = 7992.8568
$ws.Range("M136").Value = -5442.8568

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 290
$ws.Range("I7").Value = 131.16667
$ws.Range("J7").Value = 426.14285
$ws.Range("K7").Value = 393.50001
$ws.Range("L7").Value = 1278.42855
$ws.Range("M7").Value = -281.50001
$ws.Range("N7").Value = -1502.42855
$ws.Range("H36").Value = 287.5
$ws.Range("I36").Value = 245
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 735
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -566
$ws.Range("N36").Value = -1838
$ws.Range("H38").Value = 542
$ws.Range("J38").Value = 1165
$ws.Range("L38").Value = 3495
$ws.Range("N38").Value = -4189
$ws.Range("H49").Value = 1026.5
$ws.Range("I49").Value = 303
$ws.Range("J49").Value = 1750
$ws.Range("K49").Value = 909
$ws.Range("L49").Value = 5250
$ws.Range("M49").Value = -5562
$ws.Range("H50").Value = 1191.5
$ws.Range("I50").Value = 487.25
$ws.Range("K50").Value = 1461.75
$ws.Range("M50").Value = -980.75
$ws.Range("H53").Value = 1191.5
$ws.Range("I53").Value = 487.25
$ws.Range("K53").Value = 1461.75
$ws.Range("M53").Value = -980.75
$ws.Range("H132").Value = 1110.875
$ws.Range("J132").Value = 1316.3334
$ws.Range("L132").Value = 11847.0006
$ws.Range("N132").Value = -16907.0006
$ws.Range("H137").Value = 2672
$ws.Range("I137").Value = 3494
$ws.Range("J137").Value = 1850
$ws.Range("K137").Value = 10482
$ws.Range("L137").Value = 5550
$ws.Range("M137").Value = -5382
$ws.Range("N137").Value = -15750

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1950.4117
$ws.Range("I132").Value = 2029.6
$ws.Range("K132").Value = 6088.799999999999
$ws.Range("M132").Value = -3558.799999999999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4612
$ws.Range("I7").Value = 4482.6665
$ws.Range("K7").Value = 4482.6665
$ws.Range("M7").Value = -4370.6665
$ws.Range("H40").Value = 3500
$ws.Range("I40").Value = 2998
$ws.Range("K40").Value = 2998
$ws.Range("M40").Value = -2862
$ws.Range("H46").Value = 3831.12
$ws.Range("J46").Value = 3944.4546
$ws.Range("L46").Value = 3944.4546
$ws.Range("N46").Value = -4320.4546
$ws.Range("H122").Value = 3754.25
$ws.Range("J122").Value = 4505
$ws.Range("L122").Value = 13515
$ws.Range("N122").Value = -18415
$ws.Range("H126").Value = 4612
$ws.Range("I126").Value = 4482.6665
$ws.Range("K126").Value = 13447.9995
$ws.Range("M126").Value = -10977.9995
$ws.Range("H136").Value = 2456.0908
$ws.Range("I136").Value = 2279.6667
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 6839.000100000001
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -4289.000100000001
$ws.Range("N136").Value = -14850

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4633
$ws.Range("I62").Value = 4685.2856
$ws.Range("J62").Value = 4450
$ws.Range("K62").Value = 4685.2856
$ws.Range("L62").Value = 4450
$ws.Range("M62").Value = -4061.2856
$ws.Range("N62").Value = -5698
$ws.Range("H65").Value = 4633
$ws.Range("I65").Value = 4685.2856
$ws.Range("J65").Value = 4450
$ws.Range("K65").Value = 23426.428
$ws.Range("L65").Value = 22250
$ws.Range("M65").Value = -20306.428
$ws.Range("N65").Value = -28490
$ws.Range("H69").Value = 23253.2
$ws.Range("J69").Value = 23253.2
$ws.Range("L69").Value = 23253.2
$ws.Range("N69").Value = -24751.2
$ws.Range("H72").Value = 23253.2
$ws.Range("J72").Value = 23253.2
$ws.Range("L72").Value = 69759.60000000001
$ws.Range("N72").Value = -77247.60000000001
$ws.Range("H122").Value = 2798.4
$ws.Range("J122").Value = 2997
$ws.Range("L122").Value = 8991
$ws.Range("N122").Value = -13891
$ws.Range("H136").Value = 5920.385
$ws.Range("I136").Value = 5633.1816
$ws.Range("K136").Value = 16899.5448
$ws.Range("M136").Value = -14349.5448
